# Threat Alert Report update - 2026-01-18 01:00
#
# A new threat record for 13-FEB-26 (EgyptAir MS-812) is inserted ahead of
# the existing 27-FEB-26 records, pushing every row below it down by one.
# The last existing row (20-MAR-26, EgyptAir MS-812) is duplicated to the
# new row 7 to preserve it after the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Create row 7 by cloning row 6's formatting, so the new last row gets
#     the same borders/fill/font as the rest of the data rows. Values get
#     overwritten below.
$ws.Range("A6:K6").Copy()
$ws.Range("A7:K7").PasteSpecial(-4122)   # xlPasteFormats

# Force column A to text for the rows being (re)written so the "DD-MON-YY"
# strings round-trip as literal text instead of being parsed into date
# serials (matches the source file, where these are inline strings).
$ws.Range("A3:A7").NumberFormat = "@"

# --- Row 3: new record -> 13-FEB-26 / EgyptAir MS-812
$ws.Range("A3").Value = "13-FEB-26"
$ws.Range("B3").Value = "SM-328"
$ws.Range("C3").Value = "EgyptAir MS-812"
$ws.Range("D3").Value = 467
$ws.Range("E3").Value = 546
$ws.Range("F3").Value = -79
$ws.Range("G3").Value = 46
$ws.Range("H3").Value = 30
$ws.Range("I3").Value = -16
$ws.Range("J3").Value = "LOW THREAT"
$ws.Range("K3").Value = "SAR"

# --- Row 4: Air Arabia Egypt E5-590 / 27-FEB-26
$ws.Range("A4").Value = "27-FEB-26"
$ws.Range("B4").Value = "SM-328"
$ws.Range("C4").Value = "Air Arabia Egypt E5-590"
$ws.Range("D4").Value = 355
$ws.Range("E4").Value = 602
$ws.Range("F4").Value = -247
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "LOW THREAT"
$ws.Range("K4").Value = "SAR"

# --- Row 5: EgyptAir MS-812 / 27-FEB-26
$ws.Range("A5").Value = "27-FEB-26"
$ws.Range("B5").Value = "SM-328"
$ws.Range("C5").Value = "EgyptAir MS-812"
$ws.Range("D5").Value = 601
$ws.Range("E5").Value = 602
$ws.Range("F5").Value = -1
$ws.Range("G5").Value = 46
$ws.Range("H5").Value = 30
$ws.Range("I5").Value = -16
$ws.Range("J5").Value = "LOW THREAT"
$ws.Range("K5").Value = "SAR"

# --- Row 6: EgyptAir MS-812 / 13-MAR-26
$ws.Range("A6").Value = "13-MAR-26"
$ws.Range("B6").Value = "SM-328"
$ws.Range("C6").Value = "EgyptAir MS-812"
$ws.Range("D6").Value = 739
$ws.Range("E6").Value = 826
$ws.Range("F6").Value = -87
$ws.Range("G6").Value = 46
$ws.Range("H6").Value = 30
$ws.Range("I6").Value = -16
$ws.Range("J6").Value = "LOW THREAT"
$ws.Range("K6").Value = "SAR"

# --- Row 7 (new): EgyptAir MS-812 / 20-MAR-26
$ws.Range("A7").Value = "20-MAR-26"
$ws.Range("B7").Value = "SM-328"
$ws.Range("C7").Value = "EgyptAir MS-812"
$ws.Range("D7").Value = 739
$ws.Range("E7").Value = 826
$ws.Range("F7").Value = -87
$ws.Range("G7").Value = 46
$ws.Range("H7").Value = 30
$ws.Range("I7").Value = -16
$ws.Range("J7").Value = "LOW THREAT"
$ws.Range("K7").Value = "SAR"

# --- The NumberFormat="@" coercion above leaves each date cell's format
#     looking like a "text" override; restore the original look (centered,
#     bordered, no special number format) by re-pasting just the format
#     from a cell in the same column that already carries the correct
#     style (row 2 is untouched by this edit).
$ws.Range("A2").Copy()
$ws.Range("A3:A7").PasteSpecial(-4122)   # xlPasteFormats
